$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -81.8366
$ws.Range("B2").Value = -81.8179

$ws.Range("A3").Value = 25.7693
$ws.Range("B3").Value = 25.7865

$ws.Range("A4").Value = -81.4196
$ws.Range("B4").Value = -81.4382

$ws.Range("A5").Value = 26.1331
$ws.Range("B5").Value = 26.1159
